$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark that Word leaves around "Query" heading.
$d.Bookmarks("_GoBack").Delete()

# 2. Split "Input and output can not be just sequences..." so that "can not"
#    becomes "cannot", landing in its own run (matching how Word splits runs
#    when a correction is accepted), with the surrounding text kept in two
#    other runs.
$story = $d.Content
$found = $story.Find.Execute(
    "Input and output can not be just sequences but graphs, which can model cycles, branches, etc.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $target = $d.Range($story.Start, $story.End)
    $xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Input and output </w:t></w:r><w:r><w:t>cannot</w:t></w:r><w:r><w:t xml:space="preserve"> be just sequences but graphs, which can model cycles, branches, etc.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xmlFrag)
}
